# ---------------------------------------------------------------------------
# format_impor_peserta_bantuan.xlsx - add "Program" sheet, rename the
# original sheet to "Peserta", refresh the sample data / notes, and apply
# matching visual styling (bold yellow header, bold green note headers,
# boxed borders, text-forced numeric-looking codes).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Rename the original sheet, insert the new "Program" sheet before it.
$peserta = $wb.Worksheets.Item(1)
$peserta.Name = "Peserta"

$program = $wb.Worksheets.Add($peserta)     # inserted immediately before "Peserta"
$program.Name = "Program"

# ---------------------------------------------------------------------------
# 2. "Program" sheet content
# ---------------------------------------------------------------------------

$labels = @(
  "id", "Nama Program", "Sasaran Program", "Keterangan", "Asal Dana",
  "Rentang Waktu (Awal)", "Rentang Waktu (Akhir)", "Status"
)
$values = @(
  "7", "Bantuan Sosial Tunai", "2", "Bantuan Kementerian Sosial Republik Indonesia",
  "Pusat", "2020-04-01", "2020-12-31", "1"
)

# Column B holds a couple of values that look numeric/date-like ("7", "2",
# "1", the two dates) - force the whole column to Text first so Excel does
# not silently convert them to numbers / serial dates.
$program.Range("B1:B8").NumberFormat = "@"

for ($i = 0; $i -lt 8; $i++) {
    $r = $i + 1
    $program.Cells.Item($r, 1).Value = $labels[$i]
    $program.Cells.Item($r, 2).Value = $values[$i]
}

$program.Range("A9").Value = "###"
$program.Range("A10").Value = "Catatan:"
$program.Range("A11").Value = "1. Sasaran : 1 = Penduduk, 2 = Keluarga, 3 = Rumah Tangga, 4 = Kelompok"
$program.Range("A12").Value = "2. Asal Dana : 1 = Pusat, 2 = Provinsi, 3 = Kab/Kota, 4 = Dana Desa, 5 = Lain-lain(Hibah)"
$program.Range("A13").Value = "3. Status : 1 = Aktif, 2 = Tidak Aktif (Status diatur otomatis sesuai rentang waktu)"

# -- formatting: id / Nama Program / ... rows (A1:B8) - left/center aligned
$program.Range("A1:B8").HorizontalAlignment = -4131   # xlLeft
$program.Range("A1:B8").VerticalAlignment = -4108     # xlCenter

# -- "###" row (A9:B9)
$program.Range("A9:B9").HorizontalAlignment = -4131
$program.Range("A9:B9").VerticalAlignment = -4108

# -- "Catatan:" row (A10:B10) - bold, merged
$program.Range("A10:B10").Font.Bold = $true
$program.Range("A10:B10").HorizontalAlignment = -4131
$program.Range("A10:B10").VerticalAlignment = -4108
$program.Range("A10:B10").Merge()

# -- note rows (A11:A13 text, B11:B13 blank) - green fill
$program.Range("A11:B13").Interior.Color = 5296274     # RGB(146,208,80) / FF92D050
$program.Range("A11:B13").VerticalAlignment = -4108
$program.Range("A11:A13").HorizontalAlignment = -4131

# column widths
$program.Columns.Item(1).ColumnWidth = 24.875
$program.Columns.Item(2).ColumnWidth = 47.875

$program.Range("B19").Select()

# ---------------------------------------------------------------------------
# 3. "Peserta" sheet content - fix the note typo, add the new 3rd note line
# ---------------------------------------------------------------------------

$peserta.Range("A7").Value = "2. Kolom Peserta (A) dan kolom NIK (C) wajib di isi, yang lain jika kosong data diambil dari data penduduk berdasarkan kolom NIK (C)"
$peserta.Range("A7:G7").Interior.ColorIndex = 0
$peserta.Range("A7:G7").Font.Bold = $false

# header row (row 1) - bold, left/center aligned, boxed border, yellow fill
$peserta.Range("A1:G1").Font.Bold = $true
$peserta.Range("A1:G1").HorizontalAlignment = -4131
$peserta.Range("A1:G1").VerticalAlignment = -4108
$peserta.Range("A1:G1").Borders.LineStyle = 1
$peserta.Range("A1:G1").Interior.Color = 65535   # RGB(255,255,0) / FFFFFF00

# data rows (2-3) - left/center aligned, boxed border
$peserta.Range("A2:G3").HorizontalAlignment = -4131
$peserta.Range("A2:G3").VerticalAlignment = -4108
$peserta.Range("A2:G3").Borders.LineStyle = 1

# "###" row (4) - top border only, vertical center
$peserta.Range("A4:G4").Borders.Item(8).LineStyle = 1
$peserta.Range("A4:G4").VerticalAlignment = -4108

# "Catatan:" row (5) - bold, green fill, vertical center
$peserta.Range("A5:G5").Font.Bold = $true
$peserta.Range("A5:G5").Interior.Color = 5296274
$peserta.Range("A5:G5").VerticalAlignment = -4108

# note rows (6-7) - green fill, vertical center
$peserta.Range("A6:G7").Interior.Color = 5296274
$peserta.Range("A6:G7").VerticalAlignment = -4108

$peserta.Range("B13").Select()
$peserta.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 4. Final workbook view bits
# ---------------------------------------------------------------------------

$peserta.Activate()

Write-Output "done"
